# Daily update at 8 AM UTC
# Adds the next day's row (row 19) to the "Wins Over Time" sheet and moves
# the "last row" date-only formatting from row 18 (now a regular data row)
# down to the new last row 19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 is no longer the last row, so it goes back to the regular
# date+time number format used by all the other non-final rows.
$ws.Range("A18").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 19.
$ws.Range("A19").Value = 45968
$ws.Range("B19").Value = 38
$ws.Range("C19").Value = 47
$ws.Range("D19").Value = 47

# Row 19 is now the last row, so it gets the date-only number format.
$ws.Range("A19").NumberFormat = "YYYY-MM-DD"
